$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column G ("K") values for rows 2-8 (replacing old Strike# values with new K values)
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 6
$ws.Range("G5").Value = 1
$ws.Range("G6").Value = 5
$ws.Range("G7").Value = 2
$ws.Range("G8").Value = 3
